$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Estoque DPF" row (row 2) PAFmin/PAFmax values: the axis labels were
# rescaled to a smaller order of magnitude with one decimal of precision
# (e.g. 4100 -> 4.1, 4300 -> 4.3) so the chart text takes less room.
$ws.Range("C2").Value = 4.0999999999999996
$ws.Range("D2").Value = 4.3

# Active cell/selection moves from E3 to D3.
$ws.Range("D3").Select()
